$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B608").Value = "1015072666"
$ws.Range("B608").Font.Underline = $true
